# Auto-update draw results: append the 2025-11-15 Pick 4 draw as a new
# row (row 60) at the bottom of the results table, extending the used
# range from A1:E59 to A1:E60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (date-like "2025-11-15") and C (numeric-like "251115") would
# otherwise be auto-detected as a date / number by Excel's type inference.
# Force them to Text first so the new row stores plain strings, matching
# the rest of the sheet (every existing cell is a literal string).
$ws.Range("A60").NumberFormat = "@"
$ws.Range("C60").NumberFormat = "@"

$ws.Range("A60").Value = "2025-11-15"
$ws.Range("B60").Value = "Pick 4"
$ws.Range("C60").Value = "251115"
$ws.Range("D60").Value = "8-0-9-7"
$ws.Range("E60").Value = "2025-11-15T21:35:31.864+04:00"
